$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("startup")

# The "FilesTab" Cypher query (cell B4) dropped the `File Type` and `Breed`
# columns from its RETURN clause.
$newQuery = @'
MATCH (f:file)-->(parent)
WITH DISTINCT f, parent
MATCH (f)-[*]->(c:case)<--(demo:demographic)
 MATCH (s:study)<-[*]-(c)<--(diag:diagnosis)
WHERE s.clinical_study_designation IN ['COTC007B']
WITH DISTINCT f, parent, c, demo, diag, s
RETURN coalesce(f.file_name, '') AS `File Name`, 
        coalesce(labels(parent)[0], '') AS `Association`,
        coalesce(f.file_description, '') AS `Description`,
        coalesce(f.file_format, '') AS `Format`,
        coalesce(f.file_size, '') AS `Size`,
        coalesce(c.case_id, '') AS `Case ID`, 
        coalesce(diag.disease_term,'') AS Diagnosis , 
        coalesce(s.clinical_study_designation,'') AS `Study Code`
'@

$ws.Range("B4").Value2 = $newQuery

# Shorter text reflows to a shorter row.
$ws.Rows.Item(4).RowHeight = 203

# The active selection moved from B3 to B4 (with the view scrolled down a row).
$ws.Range("B4").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 4
$win.ScrollColumn = 1
